$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (target stored width = ColumnWidth + 5/6) ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# --- Update data rows 2-5 with new values (dataset refresh / custom accuracy) ---
# Row 2
$ws.Range("A2").Value = 45073.50694444445
$ws.Range("B2").Value = 5.885
$ws.Range("C2").Value = 1.942
$ws.Range("D2").Value = 1.363
$ws.Range("E2").Value = 7.767
$ws.Range("F2").Value = 3.737
$ws.Range("G2").Value = 0.96
$ws.Range("H2").Value = 6.281
$ws.Range("I2").Value = 2.168
$ws.Range("J2").Value = 0.758
$ws.Range("K2").Value = 1.218
$ws.Range("L2").Value = 2.359
$ws.Range("M2").Value = 5.493
$ws.Range("N2").Value = 0.668
$ws.Range("O2").Value = 0.884
$ws.Range("P2").Value = 2.495
$ws.Range("Q2").Value = 1.484
$ws.Range("R2").Value = 1.256
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 25.797
$ws.Range("U2").Value = 5.514
$ws.Range("V2").Value = 3.303
$ws.Range("W2").Value = 3.938
$ws.Range("X2").Value = 2.803
$ws.Range("Y2").Value = 0.246
$ws.Range("Z2").Value = 1.46
$ws.Range("AA2").Value = 1.634
$ws.Range("AB2").Value = 0.588
$ws.Range("AC2").Value = 2.048
$ws.Range("AD2").Value = 3.034
$ws.Range("AE2").Value = 0.186
$ws.Range("AF2").Value = 2.547
$ws.Range("AG2").Value = 0.594
$ws.Range("AH2").Value = 1.826

# Row 3
$ws.Range("A3").Value = 45073.51388888889
$ws.Range("B3").Value = 13.211
$ws.Range("C3").Value = 9.097
$ws.Range("D3").Value = 0.918
$ws.Range("E3").Value = 26.829
$ws.Range("F3").Value = 21.115
$ws.Range("G3").Value = 8.799
$ws.Range("H3").Value = 31.234
$ws.Range("I3").Value = 13.993
$ws.Range("J3").Value = 6.095
$ws.Range("K3").Value = 9.198
$ws.Range("L3").Value = 10.468
$ws.Range("M3").Value = 12.248
$ws.Range("N3").Value = 3.01
$ws.Range("O3").Value = 8.631
$ws.Range("P3").Value = 13.28
$ws.Range("Q3").Value = 7.666
$ws.Range("R3").Value = 0.693
$ws.Range("S3").Value = 0.213
$ws.Range("T3").Value = 135.292
$ws.Range("U3").Value = 26.033
$ws.Range("V3").Value = 9.113
$ws.Range("W3").Value = 17.871
$ws.Range("X3").Value = 9.644
$ws.Range("Y3").Value = 1.223
$ws.Range("Z3").Value = 15.398
$ws.Range("AA3").Value = 7.628
$ws.Range("AB3").Value = 6.169
$ws.Range("AC3").Value = 8.035
$ws.Range("AD3").Value = 11.254
$ws.Range("AE3").Value = 0.173
$ws.Range("AF3").Value = 26.855
$ws.Range("AG3").Value = 4.641
$ws.Range("AH3").Value = 10.489

# Row 4
$ws.Range("A4").Value = 45073.52083333334
$ws.Range("B4").Value = 19.005
$ws.Range("C4").Value = 13.81
$ws.Range("D4").Value = 0.925
$ws.Range("E4").Value = 40.173
$ws.Range("F4").Value = 32.553
$ws.Range("G4").Value = 13.931
$ws.Range("H4").Value = 54.135
$ws.Range("I4").Value = 21.769
$ws.Range("J4").Value = 9.627
$ws.Range("K4").Value = 14.451
$ws.Range("L4").Value = 15.931
$ws.Range("M4").Value = 17.59
$ws.Range("N4").Value = 4.586
$ws.Range("O4").Value = 13.745
$ws.Range("P4").Value = 20.381
$ws.Range("Q4").Value = 11.773
$ws.Range("R4").Value = 0.553
$ws.Range("S4").Value = 0.414
$ws.Range("T4").Value = 209.88
$ws.Range("U4").Value = 39.892
$ws.Range("V4").Value = 13.451
$ws.Range("W4").Value = 27.204
$ws.Range("X4").Value = 14.414
$ws.Range("Y4").Value = 1.869
$ws.Range("Z4").Value = 26.403
$ws.Range("AA4").Value = 11.634
$ws.Range("AB4").Value = 9.844
$ws.Range("AC4").Value = 12.13
$ws.Range("AD4").Value = 16.912
$ws.Range("AE4").Value = 0.141
$ws.Range("AF4").Value = 48.439
$ws.Range("AG4").Value = 7.312
$ws.Range("AH4").Value = 16.255

# Row 5
$ws.Range("A5").Value = 45073.52777777778
$ws.Range("B5").Value = 23.6
$ws.Range("C5").Value = 17.42
$ws.Range("D5").Value = 0.99
$ws.Range("E5").Value = 50.51
$ws.Range("F5").Value = 41.28
$ws.Range("G5").Value = 17.82
$ws.Range("H5").Value = 70.31
$ws.Range("I5").Value = 27.69
$ws.Range("J5").Value = 12.31
$ws.Range("K5").Value = 18.44
$ws.Range("L5").Value = 20.13
$ws.Range("M5").Value = 21.83
$ws.Range("N5").Value = 5.8
$ws.Range("O5").Value = 17.63
$ws.Range("P5").Value = 25.79
$ws.Range("Q5").Value = 14.91
$ws.Range("R5").Value = 0.5
$ws.Range("S5").Value = 0.57
$ws.Range("T5").Value = 267.19
$ws.Range("U5").Value = 50.45
$ws.Range("V5").Value = 16.85
$ws.Range("W5").Value = 34.32
$ws.Range("X5").Value = 18.1
$ws.Range("Y5").Value = 2.36
$ws.Range("Z5").Value = 34.06
$ws.Range("AA5").Value = 14.71
$ws.Range("AB5").Value = 12.63
$ws.Range("AC5").Value = 15.28
$ws.Range("AD5").Value = 21.28
$ws.Range("AE5").Value = 0.12
$ws.Range("AF5").Value = 63.27
$ws.Range("AG5").Value = 9.34
$ws.Range("AH5").Value = 20.65

# --- Remove row 6 (dataset now has one fewer timestep) ---
$ws.Rows.Item(6).Delete()

